$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.113818
$ws.Range("H2").Value = 0.341454
$ws.Range("I2").Value = 0.0003230180320166274
$ws.Range("J2").Value = 0.0003230180320166274
$ws.Range("M2").Value = 0.6068319999999999
$ws.Range("N2").Value = 1.820496
$ws.Range("O2").Value = 0.03392274820144286
$ws.Range("P2").Value = 0.03392274820144286
$ws.Range("Q2").Value = 0.06906840457599998
$ws.Range("R2").Value = 0.6216156411839999
$ws.Range("S2").Value = 0.00001095765936462566
$ws.Range("T2").Value = 0.00001095765936462566
$ws.Range("G3").Value = 0.113818
$ws.Range("H3").Value = 0.341454
$ws.Range("I3").Value = 0.0003230180320166274
$ws.Range("J3").Value = 0.0003230180320166274
$ws.Range("O3").Value = 0.4504903529585388
$ws.Range("P3").Value = 0.4504903529585388
$ws.Range("Q3").Value = 0.9172207915159999
$ws.Range("R3").Value = 8.254987123644
$ws.Range("S3").Value = 0.000145516507255143
$ws.Range("T3").Value = 0.000145516507255143
$ws.Range("G4").Value = 0.113818
$ws.Range("H4").Value = 0.341454
$ws.Range("I4").Value = 0.0003230180320166274
$ws.Range("J4").Value = 0.0003230180320166274
$ws.Range("M4").Value = 9.223151
$ws.Range("N4").Value = 27.669453
$ws.Range("O4").Value = 0.5155868988400183
$ws.Range("P4").Value = 0.5155868988400183
$ws.Range("Q4").Value = 1.049760600518
$ws.Range("R4").Value = 9.447845404661999
$ws.Range("S4").Value = 0.0001665438653968586
$ws.Range("T4").Value = 0.0001665438653968586
$ws.Range("I5").Value = 0.9904058666599795
$ws.Range("J5").Value = 0.9904058666599794
$ws.Range("M5").Value = 0.6068319999999999
$ws.Range("N5").Value = 1.820496
$ws.Range("O5").Value = 0.03392274820144286
$ws.Range("P5").Value = 0.03392274820144286
$ws.Range("Q5").Value = 211.770694861376
$ws.Range("R5").Value = 1905.936253752384
$ws.Range("S5").Value = 0.03359728883193829
$ws.Range("T5").Value = 0.03359728883193828
$ws.Range("I6").Value = 0.9904058666599795
$ws.Range("J6").Value = 0.9904058666599794
$ws.Range("O6").Value = 0.4504903529585388
$ws.Range("P6").Value = 0.4504903529585388
$ws.Range("S6").Value = 0.4461682884438617
$ws.Range("T6").Value = 0.4461682884438616
$ws.Range("I7").Value = 0.9904058666599795
$ws.Range("J7").Value = 0.9904058666599794
$ws.Range("M7").Value = 9.223151
$ws.Range("N7").Value = 27.669453
$ws.Range("O7").Value = 0.5155868988400183
$ws.Range("P7").Value = 0.5155868988400183
$ws.Range("Q7").Value = 3218.671882961668
$ws.Range("R7").Value = 28968.04694665501
$ws.Range("S7").Value = 0.5106402893841795
$ws.Range("T7").Value = 0.5106402893841795
$ws.Range("G8").Value = 3.266752
$ws.Range("H8").Value = 9.800255999999999
$ws.Range("I8").Value = 0.009271115308003845
$ws.Range("J8").Value = 0.009271115308003843
$ws.Range("M8").Value = 0.6068319999999999
$ws.Range("N8").Value = 1.820496
$ws.Range("O8").Value = 0.03392274820144286
$ws.Range("P8").Value = 0.03392274820144286
$ws.Range("Q8").Value = 1.982369649664
$ws.Range("R8").Value = 17.841326846976
$ws.Range("S8").Value = 0.0003145017101399568
$ws.Range("T8").Value = 0.0003145017101399568
$ws.Range("G9").Value = 3.266752
$ws.Range("H9").Value = 9.800255999999999
$ws.Range("I9").Value = 0.009271115308003845
$ws.Range("J9").Value = 0.009271115308003843
$ws.Range("O9").Value = 0.4504903529585388
$ws.Range("P9").Value = 0.4504903529585388
$ws.Range("Q9").Value = 26.325650205824
$ws.Range("R9").Value = 236.930851852416
$ws.Range("S9").Value = 0.004176548007421964
$ws.Range("T9").Value = 0.004176548007421963
$ws.Range("G10").Value = 3.266752
$ws.Range("H10").Value = 9.800255999999999
$ws.Range("I10").Value = 0.009271115308003845
$ws.Range("J10").Value = 0.009271115308003843
$ws.Range("M10").Value = 9.223151
$ws.Range("N10").Value = 27.669453
$ws.Range("O10").Value = 0.5155868988400183
$ws.Range("P10").Value = 0.5155868988400183
$ws.Range("Q10").Value = 30.129746975552
$ws.Range("R10").Value = 271.167722779968
$ws.Range("S10").Value = 0.004780065590441924
$ws.Range("T10").Value = 0.004780065590441922
